# Update sysC_lca.xlsx ("LCA" sheet) to reflect the new biosteam results.
#
# The transportation sub-table (previously two line items "C3"/"C4" summed
# into one "Total" row, repeated twice, then an overall "Total" row) is
# restructured into two independent item/"Total" pairs (rows 30-33) followed
# directly by the "Sum"/"All" row (row 34) - i.e. the three extra rows that
# used to sit between them are removed. Every row below shifts up by three
# to close the gap, and the numbers in the Stream and Other sub-tables are
# refreshed to the new biosteam run's results.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("LCA")

# --- 1. Split the old A30:A36 merge so A32 becomes independently editable -
$ws.Range("A30:A36").UnMerge()

# --- 2. Remove the 3 now-redundant rows (old rows 34, 35, 36) -------------
$ws.Range("A34:A36").EntireRow.Delete()

# --- 3. Refresh the transportation sub-table (rows 30-33) -----------------
$ws.Range("A30").Value = "item31 [tonne*km]"
$ws.Range("D30").Value = 1

$ws.Range("B31").Value = "Total"
$ws.Range("C31").Value = 9329531.894329507
$ws.Range("D31").Value = 1
$ws.Range("E31").Value = 1809929.187499924
$ws.Range("F31").Value = 0.1906759318900599

$ws.Range("A32").Value = "item32 [tonne*km]"
$ws.Range("B32").Value = "C4"
$ws.Range("C32").Value = 39599201.80504864
$ws.Range("E32").Value = 7682245.150179437
$ws.Range("F32").Value = 0.8093240681099401

$ws.Range("B33").Value = "Total"
$ws.Range("C33").Value = 39599201.80504864
$ws.Range("D33").Value = 1
$ws.Range("E33").Value = 7682245.150179437
$ws.Range("F33").Value = 0.8093240681099401

# Re-merge as two separate A:A pairs instead of one 4-row block
$ws.Range("A30:A31").Merge()
$ws.Range("A32:A33").Merge()

# Merge() redistributes the header border style across the new merged
# blocks (top/bottom split); restore the plain header style used
# throughout the sheet by copying it from an unaffected header cell.
$ws.Range("A29").Copy()
$ws.Range("A30:A33").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# --- 4. Refresh the Stream sub-table (rows 38-46, after the row shift) ----
$ws.Range("B38").Value = 1737293.784823114
$ws.Range("C38").Value = 48644225.97504719
$ws.Range("D38").Value = 26.64626608726223

$ws.Range("D39").Value = 1.219904601981638
$ws.Range("D40").Value = -1.751484086509394
$ws.Range("D41").Value = -20.57714712907902
$ws.Range("D42").Value = -0.6792282764243533
$ws.Range("D43").Value = -0.6279452047606762
$ws.Range("D44").Value = -1.8313914137129
$ws.Range("D45").Value = -1.398974578757533

$ws.Range("C46").Value = 1825555.063352786

# --- 5. Rename the "Other" sub-table's single item (row 50) ---------------
$ws.Range("A50").Value = "E_item [kWh]"
